$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Update row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.4
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Update row 4
$ws.Range("B4").Value = 0.3
$ws.Range("C4").Value = 0.4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Remove rows 5 and 6 (old nodes 4 and 5) which no longer exist in the target
$ws.Range("A5:E6").ClearContents()

# Apply the underline style used in D11 to D10 as well (new empty styled cell)
$ws.Range("D10").Font.Underline = $true

# Update the active selection to match the target view
$ws.Range("F7").Select()
